$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$lo = $ws.ListObjects.Item(1)

# Insert a new row at sheet row 71 (shifts rows 71..137 down to 72..138),
# extending the Table1 range from A8:K137 to A8:K138.
$ws.Rows.Item(71).Insert(-4121, 0)
$lo.Resize($ws.Range("A8:K138"))

# Copy number formats/styles from the row below (the row that used to be
# row 71, now shifted to row 72) into the newly blank row 71, matching
# the rest of the table's per-column styling (only columns A:K).
$ws.Range("A72:K72").Copy()
$ws.Range("A71:K71").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Restore the calculated-column formula in the new row's "EARNED " column
# (G), which PasteSpecial(Formats) does not bring along.
$ws.Range("G71").Formula = "=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"""",Table1[[#This Row],[EARNED]])"

# --- Row 69: record 1.25 days SL earned ---
$ws.Range("C69").Value2 = 1.25

# --- Row 70: SL(3-0-0) used 8/3,4,7/2023 ---
$ws.Range("B70").Value2 = "SL(3-0-0)"
$ws.Range("C70").Value2 = 1.25
$ws.Range("H70").Value2 = 3
$ws.Range("K70").Value2 = "8/3,4,7/2023"

# --- Row 71 (new row): SL(2-0-0) used 8/24,25/2023 ---
$ws.Range("B71").Value2 = "SL(2-0-0)"
$ws.Range("H71").Value2 = 2
$ws.Range("K71").Value2 = "8/24,25/2023"
